# Update "Latest HO Xliff Generate Date" / Handoff / Handback datetime
# values to reflect the freshly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for ce66918e row (row 3)
$wsOverview.Range("G3").Value = "2016-09-06 09:01:47"

# zh-cn sheet: ce66918e row (row 3)
# Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn.Range("H3").Value = "2016-09-06 09:01:41"
$wsZhCn.Range("K3").Value = "2016-09-06 09:02:05"

# de-de sheet: ce66918e row (row 3)
# Correspond Handoff Datetime (H3) mirrors Overview's Latest HO Xliff Generate Date
$wsDeDe.Range("H3").Value = "2016-09-06 09:01:47"
# Correspond Handback DateTime (K3)
$wsDeDe.Range("K3").Value = "2016-09-06 09:02:28"
